$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply the date/time number format to column A data rows (A2:A13).
# First set with a lowercase format code (registers numFmtId 164), then
# switch to the uppercase variant actually used (registers numFmtId 165),
# matching the target style sheet which retains both entries.
$ws.Cells.Item(2,1).NumberFormat = "yyyy-mm-dd h:mm:ss"
$ws.Cells.Item(2,1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
for ($r = 3; $r -le 13; $r++) {
  $ws.Cells.Item($r, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
}

# New "Trening" column header - copy the header style from E1 (bold font,
# border, centered alignment) onto F1, then set its text.
$ws.Range("E1").Copy($ws.Range("F1"))
$ws.Cells.Item(1,6).Value2 = "Trening"

$ws.Cells.Item(2,1).Value2 = 45685.64971921296
$ws.Cells.Item(2,2).Value2 = 1250.7
$ws.Cells.Item(2,3).Value2 = 14.29
$ws.Cells.Item(2,4).Value2 = 3.63045460837228
$ws.Cells.Item(2,5).Value2 = "10-15"
$ws.Cells.Item(2,6).Value2 = "Duża Gra"

$ws.Cells.Item(3,1).Value2 = 45685.65116481482
$ws.Cells.Item(3,2).Value2 = 1375.6
$ws.Cells.Item(3,3).Value2 = 13.7
$ws.Cells.Item(3,4).Value2 = 3.8644768510546
$ws.Cells.Item(3,5).Value2 = "10-15"
$ws.Cells.Item(3,6).Value2 = "Duża Gra"

$ws.Cells.Item(4,1).Value2 = 45685.65508958334
$ws.Cells.Item(4,2).Value2 = 1714.7
$ws.Cells.Item(4,3).Value2 = 14.23
$ws.Cells.Item(4,4).Value2 = 3.947716849190849
$ws.Cells.Item(4,5).Value2 = "10-15"
$ws.Cells.Item(4,6).Value2 = "Duża Gra"

$ws.Cells.Item(5,1).Value2 = 45685.64971574074
$ws.Cells.Item(5,2).Value2 = 1250.4
$ws.Cells.Item(5,3).Value2 = 9.960000000000001
$ws.Cells.Item(5,4).Value2 = 3.123471191951207
$ws.Cells.Item(5,5).Value2 = "5-10"
$ws.Cells.Item(5,6).Value2 = "Duża Gra"

$ws.Cells.Item(6,1).Value2 = 45685.65116134259
$ws.Cells.Item(6,2).Value2 = 1375.3
$ws.Cells.Item(6,3).Value2 = 9.17
$ws.Cells.Item(6,4).Value2 = 2.955209919384548
$ws.Cells.Item(6,5).Value2 = "5-10"
$ws.Cells.Item(6,6).Value2 = "Duża Gra"

$ws.Cells.Item(7,1).Value2 = 45685.6664125
$ws.Cells.Item(7,2).Value2 = 2693
$ws.Cells.Item(7,3).Value2 = 8.699999999999999
$ws.Cells.Item(7,4).Value2 = 3.159672907420568
$ws.Cells.Item(7,5).Value2 = "5-10"
$ws.Cells.Item(7,6).Value2 = "Duża Gra"

$ws.Cells.Item(8,1).Value2 = 45685.67503402778
$ws.Cells.Item(8,2).Value2 = 3437.9
$ws.Cells.Item(8,3).Value2 = 14.58
$ws.Cells.Item(8,4).Value2 = 3.400056259972708
$ws.Cells.Item(8,5).Value2 = "10-15"
$ws.Cells.Item(8,6).Value2 = "Mała Gra"

$ws.Cells.Item(9,1).Value2 = 45685.67546805555
$ws.Cells.Item(9,2).Value2 = 3475.4
$ws.Cells.Item(9,3).Value2 = 13
$ws.Cells.Item(9,4).Value2 = 3.57370798928397
$ws.Cells.Item(9,5).Value2 = "10-15"
$ws.Cells.Item(9,6).Value2 = "Mała Gra"

$ws.Cells.Item(10,1).Value2 = 45685.67800856481
$ws.Cells.Item(10,2).Value2 = 3694.9
$ws.Cells.Item(10,3).Value2 = 11.46
$ws.Cells.Item(10,4).Value2 = 3.421598468508037
$ws.Cells.Item(10,5).Value2 = "10-15"
$ws.Cells.Item(10,6).Value2 = "Mała Gra"

$ws.Cells.Item(11,1).Value2 = 45685.67800625
$ws.Cells.Item(11,2).Value2 = 3694.7
$ws.Cells.Item(11,3).Value2 = 8.279999999999999
$ws.Cells.Item(11,4).Value2 = 2.99589272907802
$ws.Cells.Item(11,5).Value2 = "5-10"
$ws.Cells.Item(11,6).Value2 = "Mała Gra"

$ws.Cells.Item(12,1).Value2 = 45685.68072037037
$ws.Cells.Item(12,2).Value2 = 3929.2
$ws.Cells.Item(12,3).Value2 = 9.720000000000001
$ws.Cells.Item(12,4).Value2 = 2.853019612176079
$ws.Cells.Item(12,5).Value2 = "5-10"
$ws.Cells.Item(12,6).Value2 = "Mała Gra"

$ws.Cells.Item(13,1).Value2 = 45685.68168564815
$ws.Cells.Item(13,2).Value2 = 4012.6
$ws.Cells.Item(13,3).Value2 = 9.789999999999999
$ws.Cells.Item(13,4).Value2 = 3.069698538099017
$ws.Cells.Item(13,5).Value2 = "5-10"
$ws.Cells.Item(13,6).Value2 = "Mała Gra"

